# Event_Timing.xlsx edit
# - Refresh the workbook/window view state (book views block gets the full
#   set of attributes Excel normally writes out, not just activeTab).
# - Corrected D-Station derived measurements in row 4 (G4, H4, I4, Q4) now
#   that Strain / Strain Rate files are located automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$w  = $excel.ActiveWindow

# --- Window / book view state -------------------------------------------
$w.DisplayHorizontalScrollBar = $true
$w.DisplayVerticalScrollBar   = $true
$w.TabRatio                   = 600
$w.AutoFilterDateGrouping     = $true
$w.WindowState                = -4143

# --- Updated measurements (tresparametros -> D-Station) ------------------
$ws.Range("G4").Value = 102
$ws.Range("H4").Value = 1140
$ws.Range("I4").Value = 1037
$ws.Range("Q4").Value = 700
